$wb = $excel.ActiveWorkbook

# --- Sheet ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H64").Value = 28573848
$ws.Range("I64").Value = 55557530
$ws.Range("J64").Value = 2894.4707
$ws.Range("K64").Value = 55557530
$ws.Range("L64").Value = 2894.4707
$ws.Range("M64").Value = -55557282
$ws.Range("N64").Value = -3390.4707
$ws.Range("H67").Value = 28573848
$ws.Range("I67").Value = 55557530
$ws.Range("J67").Value = 2894.4707
$ws.Range("K67").Value = 55557530
$ws.Range("L67").Value = 2894.4707
$ws.Range("M67").Value = -55556672
$ws.Range("N67").Value = -4610.4707
$ws.Range("H69").Value = 3742.457
$ws.Range("I69").Value = 3754.5454
$ws.Range("J69").Value = 3736.9167
$ws.Range("K69").Value = 11263.6362
$ws.Range("L69").Value = 11210.7501
$ws.Range("M69").Value = -10389.6362
$ws.Range("N69").Value = -12958.7501
$ws.Range("H72").Value = 3742.457
$ws.Range("I72").Value = 3754.5454
$ws.Range("J72").Value = 3736.9167
$ws.Range("K72").Value = 33790.9086
$ws.Range("L72").Value = 33632.2503
$ws.Range("M72").Value = -29422.9086
$ws.Range("N72").Value = -42368.2503
$ws.Range("H74").Value = 2848.7727
$ws.Range("I74").Value = 2691.75
$ws.Range("J74").Value = 3037.2
$ws.Range("K74").Value = 2691.75
$ws.Range("L74").Value = 3037.2
$ws.Range("M74").Value = -1755.75
$ws.Range("N74").Value = -4909.2
$ws.Range("H77").Value = 2848.7727
$ws.Range("I77").Value = 2691.75
$ws.Range("J77").Value = 3037.2
$ws.Range("K77").Value = 13458.75
$ws.Range("L77").Value = 15186
$ws.Range("M77").Value = -8778.75
$ws.Range("N77").Value = -24546
$ws.Range("H86").Value = 1644.9231
$ws.Range("I86").Value = 1744.5883
$ws.Range("J86").Value = 1456.6666
$ws.Range("K86").Value = 1744.5883
$ws.Range("L86").Value = 1456.6666
$ws.Range("M86").Value = -621.5882999999999
$ws.Range("N86").Value = -3702.6666
$ws.Range("H89").Value = 1644.9231
$ws.Range("I89").Value = 1744.5883
$ws.Range("J89").Value = 1456.6666
$ws.Range("K89").Value = 8722.941499999999
$ws.Range("L89").Value = 7283.333000000001
$ws.Range("M89").Value = -3106.941499999999
$ws.Range("N89").Value = -18515.333

# --- Sheet BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 4811.222
$ws.Range("I86").Value = 3767
$ws.Range("K86").Value = 3767
$ws.Range("M86").Value = -2644
$ws.Range("H89").Value = 4811.222
$ws.Range("I89").Value = 3767
$ws.Range("K89").Value = 18835
$ws.Range("M89").Value = -13219

# --- Sheet CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 17682.098
$ws.Range("I31").Value = 1061.5358
$ws.Range("J31").Value = 28258.818
$ws.Range("K31").Value = 1061.5358
$ws.Range("L31").Value = 28258.818
$ws.Range("M31").Value = -766.5358000000001
$ws.Range("N31").Value = -28848.818
$ws.Range("H34").Value = 17682.098
$ws.Range("I34").Value = 1061.5358
$ws.Range("J34").Value = 28258.818
$ws.Range("K34").Value = 1061.5358
$ws.Range("L34").Value = 28258.818
$ws.Range("M34").Value = -859.5358000000001
$ws.Range("N34").Value = -28662.818

# --- Sheet GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 4804.185
$ws.Range("J70").Value = 4878.6665
$ws.Range("L70").Value = 4878.6665
$ws.Range("N70").Value = -5418.6665
$ws.Range("H73").Value = 4804.185
$ws.Range("J73").Value = 4878.6665
$ws.Range("L73").Value = 4878.6665
$ws.Range("N73").Value = -6750.6665
$ws.Range("H80").Value = 2991
$ws.Range("I80").Value = 2401.6667
$ws.Range("J80").Value = 3875
$ws.Range("K80").Value = 2401.6667
$ws.Range("L80").Value = 3875
$ws.Range("M80").Value = -1403.6667
$ws.Range("N80").Value = -5871
$ws.Range("H83").Value = 2991
$ws.Range("I83").Value = 2401.6667
$ws.Range("J83").Value = 3875
$ws.Range("K83").Value = 12008.3335
$ws.Range("L83").Value = 19375
$ws.Range("M83").Value = -7016.333500000001
$ws.Range("N83").Value = -29359
$ws.Range("H126").Value = 1278
$ws.Range("I126").Value = 1212
$ws.Range("J126").Value = 1300
$ws.Range("K126").Value = 3636
$ws.Range("L126").Value = 3900
$ws.Range("M126").Value = -1166
$ws.Range("N126").Value = -8840

# --- Sheet LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2802.5
$ws.Range("I7").Value = 0
$ws.Range("J7").Value = 2802.5
$ws.Range("K7").Value = 0
$ws.Range("L7").Value = 2802.5
$ws.Range("M7").ClearContents()
$ws.Range("N7").Value = -3026.5
$ws.Range("H68").Value = 1843.25
$ws.Range("I68").Value = 1734.4286
$ws.Range("J68").Value = 2097.1667
$ws.Range("K68").Value = 1734.4286
$ws.Range("L68").Value = 2097.1667
$ws.Range("M68").Value = -985.4286
$ws.Range("N68").Value = -3595.1667
$ws.Range("H71").Value = 1843.25
$ws.Range("I71").Value = 1734.4286
$ws.Range("J71").Value = 2097.1667
$ws.Range("K71").Value = 8672.143
$ws.Range("L71").Value = 10485.8335
$ws.Range("M71").Value = -4928.143
$ws.Range("N71").Value = -17973.8335
$ws.Range("H82").Value = 1469.4615
$ws.Range("I82").Value = 1444.7778
$ws.Range("J82").Value = 1525
$ws.Range("K82").Value = 1444.7778
$ws.Range("L82").Value = 1525
$ws.Range("M82").Value = -1083.7778
$ws.Range("N82").Value = -2247
$ws.Range("H85").Value = 1469.4615
$ws.Range("I85").Value = 1444.7778
$ws.Range("J85").Value = 1525
$ws.Range("K85").Value = 1444.7778
$ws.Range("L85").Value = 1525
$ws.Range("M85").Value = -196.7778000000001
$ws.Range("N85").Value = -4021
$ws.Range("H126").Value = 2802.5
$ws.Range("I126").Value = 0
$ws.Range("J126").Value = 2802.5
$ws.Range("K126").Value = 0
$ws.Range("L126").Value = 8407.5
$ws.Range("M126").ClearContents()
$ws.Range("N126").Value = -13347.5
$ws.Range("H127").Value = 45036
$ws.Range("J127").Value = 45036
$ws.Range("L127").Value = 45036
$ws.Range("N127").Value = -54956

# --- Sheet WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 2799.6924
$ws.Range("I62").Value = 2767.3333
$ws.Range("J62").Value = 2872.5
$ws.Range("K62").Value = 2767.3333
$ws.Range("L62").Value = 2872.5
$ws.Range("M62").Value = -2143.3333
$ws.Range("N62").Value = -4120.5
$ws.Range("H65").Value = 2799.6924
$ws.Range("I65").Value = 2767.3333
$ws.Range("J65").Value = 2872.5
$ws.Range("K65").Value = 13836.6665
$ws.Range("L65").Value = 14362.5
$ws.Range("M65").Value = -10716.6665
$ws.Range("N65").Value = -20602.5
$ws.Range("H126").Value = 890
$ws.Range("I126").Value = 775
$ws.Range("J126").Value = 982
$ws.Range("K126").Value = 2325
$ws.Range("L126").Value = 2946
$ws.Range("M126").Value = 145
$ws.Range("N126").Value = -7886
$ws.Range("H141").Value = 50000
$ws.Range("J141").Value = 50000
$ws.Range("L141").Value = 50000
$ws.Range("N141").Value = -60360
